$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must stay text (matches source data which
# stores prices as text). Force text format first so Excel does not re-parse them as numbers.
$textCells = @("D5", "D6", "D8", "D14", "D16", "D21", "D24", "D25", "D26", "D28", "D29", "D31", "D33", "D34", "D35", "D39", "D43", "D45", "D46", "D47", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values (same order rows appear in the sheet).
$ws.Range("D2").Value = "60.677.95"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "2.904.59"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "528.30"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").Value = "143.48"
$ws.Range("E6").Value = "  -5.60%  "
$ws.Range("D8").Value = "0.556"
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").Value = "2.915.43"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  -4.93%  "
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "3.419.58"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").Value = "0.127"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "60.656.94"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "22.80"
$ws.Range("E16").Value = "  -4.10%  "
$ws.Range("D17").Value = "2.917.17"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").Value = "362.42"
$ws.Range("E21").Value = "  -5.91%  "
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "5.69"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "64.78"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "0.455"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "7.88"
$ws.Range("E29").Value = "  -5.77%  "
$ws.Range("D30").Value = "0.0₃0851"
$ws.Range("E30").Value = "  -9.84%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "19.81"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "151.23"
$ws.Range("E34").Value = "  -5.28%  "
$ws.Range("D35").Value = "4.37"
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("E36").Value = "  -6.31%  "
$ws.Range("E37").Value = "  -6.37%  "
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("D39").Value = "37.89"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  -4.75%  "
$ws.Range("E41").Value = "  -5.42%  "
$ws.Range("D42").Value = "2.295.68"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "0.650"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "20.44"
$ws.Range("E45").Value = "  -8.04%  "
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "5.00"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").Value = "250.39"
$ws.Range("E51").Value = "  -6.85%  "
